$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")
$ws.Range("C1").NumberFormat = "mm-dd-yy"
$ws.Range("C1").Value = (Get-Date -Year 2021 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0).Date
